$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Fitness) values according to run-length encoded ranges
$ranges = @(
    @{Start=2; End=4; Value=11388}
    @{Start=5; End=8; Value=11029}
    @{Start=9; End=10; Value=10529}
    @{Start=11; End=14; Value=10178}
    @{Start=15; End=16; Value=9779}
    @{Start=17; End=20; Value=9297}
    @{Start=21; End=27; Value=8942}
    @{Start=28; End=29; Value=8657}
    @{Start=30; End=30; Value=8365}
    @{Start=31; End=40; Value=7812}
    @{Start=41; End=43; Value=7598}
    @{Start=44; End=252; Value=7573}
)

foreach ($r in $ranges) {
    $ws.Range("C$($r.Start):C$($r.End)").Value = $r.Value
}
